$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '66.718.18'
Set-TextValue 'E2' '  -1.20%  '
Set-TextValue 'D3' '2.524.32'
Set-TextValue 'E3' '  -3.91%  '
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '584.80'
Set-TextValue 'E5' '  -1.50%  '
Set-TextValue 'D6' '171.47'
Set-TextValue 'E6' '  +2.19%  '
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'E8' '  -1.44%  '
Set-TextValue 'D9' '2.523.15'
Set-TextValue 'E9' '  -3.95%  '
Set-TextValue 'D10' '0.140'
Set-TextValue 'E10' '  -0.25%  '
Set-TextValue 'E11' '  +0.02%  '
Set-TextValue 'E12' '  -3.89%  '
Set-TextValue 'D13' '5.14'
Set-TextValue 'E13' '  -1.76%  '
Set-TextValue 'D14' '26.76'
Set-TextValue 'E14' '  -3.18%  '
Set-TextValue 'D15' '2.983.95'
Set-TextValue 'E15' '  -4.29%  '
Set-TextValue 'E16' '  -2.88%  '
Set-TextValue 'D17' '66.567.27'
Set-TextValue 'E17' '  -1.58%  '
Set-TextValue 'D18' '2.520.76'
Set-TextValue 'E18' '  -3.87%  '
Set-TextValue 'D19' '7.87'
Set-TextValue 'E19' '  -2.09%  '
Set-TextValue 'E20' '  -5.68%  '
Set-TextValue 'E21' '  -2.89%  '
Set-TextValue 'D22' '4.21'
Set-TextValue 'E22' '  -2.78%  '
Set-TextValue 'D23' '4.66'
Set-TextValue 'E23' '  -0.80%  '
Set-TextValue 'D24' '1.98'
Set-TextValue 'E24' '  +2.15%  '
Set-TextValue 'E25' '  +0.14%  '
Set-TextValue 'D26' '70.23'
Set-TextValue 'E26' '  +0.32%  '
Set-TextValue 'D27' '9.95'
Set-TextValue 'D28' '1.00'
Set-TextValue 'E28' '  +0.13%  '
Set-TextValue 'D29' '2.634.52'
Set-TextValue 'E29' '  -4.61%  '
Set-TextValue 'D30' '0.0₃0980'
Set-TextValue 'E30' '  -2.75%  '
Set-TextValue 'D31' '527.16'
Set-TextValue 'E31' '  -3.58%  '
Set-TextValue 'D32' '8.15'
Set-TextValue 'E32' '  +2.34%  '
Set-TextValue 'E33' '  -2.56%  '
Set-TextValue 'E34' '  -2.79%  '
Set-TextValue 'E35' '  -4.49%  '
Set-TextValue 'D36' '0.999'
Set-TextValue 'E36' '  -0.05%  '
Set-TextValue 'E37' '  -2.75%  '
Set-TextValue 'D38' '157.53'
Set-TextValue 'E38' '  -0.28%  '
Set-TextValue 'D39' '18.64'
Set-TextValue 'E39' '  -2.20%  '
Set-TextValue 'D40' '18.39'
Set-TextValue 'E40' '  +0.63%  '
Set-TextValue 'E42' '  -0.85%  '
Set-TextValue 'E43' '  -2.14%  '
Set-TextValue 'E45' '  +3.24%  '
Set-TextValue 'D46' '39.45'
Set-TextValue 'E46' '  -1.46%  '
Set-TextValue 'D47' '149.64'
Set-TextValue 'D48' '0.560'
Set-TextValue 'E48' '  -3.67%  '
Set-TextValue 'E49' '  -2.94%  '
Set-TextValue 'E50' '  +1.20%  '
Set-TextValue 'D51' '0.0₆0269'
Set-TextValue 'E51' '  -10.80%  '
